# "add score to gn sample"
# Add a new "Score" column (M) to the evaluation results sheet, with a
# bold/centered/bordered header matching the existing header style, and
# per-row integer scores.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell M1: "Score" ---------------------------------------------
$ws.Range("M1").Value = "Score"

# Match the look of the other header cells (bold font, centered/top
# alignment) by copying L1's formatting onto M1, then trim the border down
# to just the left/right edges (matches the diff's new borderId/cellXfs).
$ws.Range("L1").Copy()
$ws.Range("M1").PasteSpecial(-4122)

$ws.Range("M1").Borders.Item(8).LineStyle = -4142
$ws.Range("M1").Borders.Item(9).LineStyle = -4142
$ws.Range("M1").Borders.Item(7).LineStyle = 1
$ws.Range("M1").Borders.Item(7).Weight = 2
$ws.Range("M1").Borders.Item(10).LineStyle = 1
$ws.Range("M1").Borders.Item(10).Weight = 2

# --- Per-row scores (M2:M21) ----------------------------------------------
$scores = 8, 11, 8, 8, 9, 9, 9, 8, 7, 10, 10, 4, 7, 7, 9, 9, 8, 11, 1, 10
for ($i = 0; $i -lt $scores.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 13).Value = $scores[$i]
}

# --- Match the author's final selection in the saved workbook ------------
$ws.Range("E29").Select() | Out-Null
